# Applies the "cryptos list" data refresh described by the commit:
# "Updated cryptos list on Mon Mar  4 13:45:53 UTC 2024 with GitHub Actions"
#
# For each changed row, Price (column D) and/or Volume(1h) (column E) are updated
# to the new scraped values. Rows 50/51 additionally swap their Coin name/Link
# (ApeXProtocol and EnergySwap traded ranking positions).
#
# Price values that look like plain numbers are written with a leading apostrophe
# so Excel keeps them as text (matching the original inline-string/text cells,
# e.g. keeping a trailing zero such as "1.00" or "0.780" instead of becoming 1 / 0.78).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "65.607.32"
$ws.Range("E2").Value = "  +5.42%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "3.514.27"
$ws.Range("E3").Value = "  +2.66%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  +0.10%  "

# Row 5: BNB
$ws.Range("D5").Value = "'418.62"
$ws.Range("E5").Value = "  +1.34%  "

# Row 6: Solana
$ws.Range("D6").Value = "'132.79"
$ws.Range("E6").Value = "  +3.22%  "

# Row 7: XRP
$ws.Range("E7").Value = "  +4.63%  "

# Row 8: USDC
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  +0.02%  "

# Row 9: Cardano
$ws.Range("D9").Value = "'0.780"
$ws.Range("E9").Value = "  +7.09%  "

# Row 10: Dogecoin
$ws.Range("D10").Value = "'0.163"
$ws.Range("E10").Value = "  +16.81%  "

# Row 11: Avalanche
$ws.Range("D11").Value = "'43.37"
$ws.Range("E11").Value = "  +1.40%  "

# Row 12: ShibaInu
$ws.Range("D12").Value = "'0.0000265"
$ws.Range("E12").Value = "  +20.06%  "

# Row 13: Polkadot
$ws.Range("D13").Value = "'10.02"
$ws.Range("E13").Value = "  +8.93%  "

# Row 14: WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "4.065.00"
$ws.Range("E14").Value = "  +2.72%  "

# Row 15: TRON
$ws.Range("E15").Value = "  +0.23%  "

# Row 16: Chainlink
$ws.Range("D16").Value = "'20.55"
$ws.Range("E16").Value = "  +0.58%  "

# Row 17: WrappedEther
$ws.Range("D17").Value = "3.489.17"
$ws.Range("E17").Value = "  +1.52%  "

# Row 18: Uniswap
$ws.Range("D18").Value = "'12.90"
$ws.Range("E18").Value = "  +1.53%  "

# Row 19: Polygon
$ws.Range("D19").Value = "'1.10"
$ws.Range("E19").Value = "  +3.06%  "

# Row 20: WrappedBTC
$ws.Range("D20").Value = "65.469.66"
$ws.Range("E20").Value = "  +5.22%  "

# Row 21: BitcoinCash
$ws.Range("D21").Value = "'454.68"
$ws.Range("E21").Value = "  -4.07%  "

# Row 22: Litecoin
$ws.Range("D22").Value = "'90.28"
$ws.Range("E22").Value = "  -1.29%  "

# Row 23: ImmutableX
$ws.Range("E23").Value = "  -0.36%  "

# Row 24: InternetComputer(DFINITY)
$ws.Range("D24").Value = "'13.31"
$ws.Range("E24").Value = "  +1.23%  "

# Row 25: PancakeSwap
$ws.Range("D25").Value = "'3.41"
$ws.Range("E25").Value = "  +3.54%  "

# Row 26: Filecoin
$ws.Range("D26").Value = "'9.95"
$ws.Range("E26").Value = "  +2.10%  "

# Row 27: EthereumClassic
$ws.Range("D27").Value = "'34.29"
$ws.Range("E27").Value = "  +2.61%  "

# Row 28: Cosmos
$ws.Range("D28").Value = "'12.64"
$ws.Range("E28").Value = "  +6.69%  "

# Row 30: RenderToken
$ws.Range("D30").Value = "'7.46"
$ws.Range("E30").Value = "  -3.98%  "

# Row 31: Hedera
$ws.Range("E31").Value = "  +5.89%  "

# Row 32: Kaspa
$ws.Range("E32").Value = "  -1.15%  "

# Row 33: InjectiveProtocol
$ws.Range("D33").Value = "'39.75"
$ws.Range("E33").Value = "  -2.80%  "

# Row 34: Dai
$ws.Range("E34").Value = "  +0.02%  "

# Row 35: OKB
$ws.Range("D35").Value = "'57.57"
$ws.Range("E35").Value = "  -0.63%  "

# Row 36: VeChain
$ws.Range("E36").Value = "  +4.17%  "

# Row 37: PEPE
$ws.Range("D37").Value = "0.0₃0741"
$ws.Range("E37").Value = "  +35.80%  "

# Row 38: Stellar
$ws.Range("E38").Value = "  +9.83%  "

# Row 39: FirstDigitalUSD
$ws.Range("D39").Value = "'0.997"
$ws.Range("E39").Value = "  -0.15%  "

# Row 40: Stacks
$ws.Range("D40").Value = "'3.07"
$ws.Range("E40").Value = "  +0.98%  "

# Row 41: NEARProtocol
$ws.Range("D41").Value = "'4.52"
$ws.Range("E41").Value = "  +4.17%  "

# Row 42: WEMIXToken
$ws.Range("E42").Value = "  +3.77%  "

# Row 43: Monero
$ws.Range("D43").Value = "'146.02"
$ws.Range("E43").Value = "  +0.08%  "

# Row 44: LidoDAOToken
$ws.Range("E44").Value = "  -0.59%  "

# Row 45: TheGraph
$ws.Range("D45").Value = "'0.312"
$ws.Range("E45").Value = "  -3.01%  "

# Row 46: ARBITRUM
$ws.Range("E46").Value = "  -2.70%  "

# Row 47: ThetaToken
$ws.Range("D47").Value = "'2.32"
$ws.Range("E47").Value = "  -1.43%  "

# Row 48: Celestia
$ws.Range("D48").Value = "'15.85"
$ws.Range("E48").Value = "  -2.80%  "

# Row 49: Cronos
$ws.Range("D49").Value = "'0.145"
$ws.Range("E49").Value = "  +3.88%  "

# Row 50: EnergySwap -> ApeXProtocol
$ws.Range("B50").Value = "ApeXProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D50").Value = "'2.57"
$ws.Range("E50").Value = "  +10.74%  "

# Row 51: ApeXProtocol -> EnergySwap
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'21.63"
$ws.Range("E51").Value = "  -2.87%  "

Write-Host "Updated cryptos list with latest price/volume data"
